$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 8400
$ws.Range("J16").Value = 8400
$ws.Range("L16").Value = 8400
$ws.Range("N16").Value = -8860
$ws.Range("H64").Value = 7245.8823
$ws.Range("I64").Value = 3909
$ws.Range("K64").Value = 3909
$ws.Range("M64").Value = -3661
$ws.Range("H67").Value = 7245.8823
$ws.Range("I67").Value = 3909
$ws.Range("K67").Value = 3909
$ws.Range("M67").Value = -3051
$ws.Range("H112").Value = 3305.9167
$ws.Range("J112").Value = 3330.8096
$ws.Range("L112").Value = 9992.4288
$ws.Range("N112").Value = -12208.4288
$ws.Range("H137").Value = 4359.0684
$ws.Range("I137").Value = 4658.5
$ws.Range("K137").Value = 13975.5
$ws.Range("M137").Value = -11425.5
$ws.Range("H138").Value = 4087.923
$ws.Range("I138").Value = 3576.3157
$ws.Range("J138").Value = 4573.95
$ws.Range("K138").Value = 10728.9471
$ws.Range("L138").Value = 13721.85
$ws.Range("M138").Value = -5588.947100000001
$ws.Range("N138").Value = -24001.85

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 5732.3335
$ws.Range("I21").Value = 1098.5
$ws.Range("K21").Value = 1098.5
$ws.Range("M21").Value = -724.5
$ws.Range("H30").Value = 1500
$ws.Range("I30").Value = 1500
$ws.Range("K30").Value = 1500
$ws.Range("M30").Value = -1350
$ws.Range("H32").Value = 7387.8804
$ws.Range("I32").Value = 7387.8804
$ws.Range("K32").Value = 7387.8804
$ws.Range("M32").Value = -7100.8804
$ws.Range("H45").Value = 3290
$ws.Range("I45").Value = 2641.6365
$ws.Range("J45").Value = 4082.4443
$ws.Range("K45").Value = 2641.6365
$ws.Range("L45").Value = 4082.4443
$ws.Range("M45").Value = -2264.6365
$ws.Range("N45").Value = -4836.4443
$ws.Range("H110").Value = 2807.5833
$ws.Range("I110").Value = 869.1
$ws.Range("K110").Value = 869.1
$ws.Range("M110").Value = 1175.9
$ws.Range("H122").Value = 2409.1035
$ws.Range("I122").Value = 2310.5
$ws.Range("K122").Value = 6931.5
$ws.Range("M122").Value = -4481.5
$ws.Range("H132").Value = 1515.75
$ws.Range("I132").Value = 1532.3871
$ws.Range("K132").Value = 4597.1613
$ws.Range("M132").Value = -2067.1613

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 657.19354
$ws.Range("I94").Value = 672.43335
$ws.Range("J94").Value = 200
$ws.Range("K94").Value = 672.43335
$ws.Range("L94").Value = 200
$ws.Range("M94").Value = -221.43335
$ws.Range("N94").Value = -1102
$ws.Range("H105").Value = 2372
$ws.Range("I105").Value = 2484.0667
$ws.Range("J105").Value = 2035.8
$ws.Range("K105").Value = 2484.0667
$ws.Range("L105").Value = 2035.8
$ws.Range("M105").Value = -737.0666999999999
$ws.Range("N105").Value = -5529.8
$ws.Range("H134").Value = 2163.2285
$ws.Range("I134").Value = 1834
$ws.Range("J134").Value = 4714.75
$ws.Range("K134").Value = 5502
$ws.Range("L134").Value = 14144.25
$ws.Range("M134").Value = -2967
$ws.Range("N134").Value = -19214.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28574696
$ws.Range("I31").Value = 50001736
$ws.Range("K31").Value = 50001736
$ws.Range("M31").Value = -50001441
$ws.Range("H34").Value = 28574696
$ws.Range("I34").Value = 50001736
$ws.Range("K34").Value = 50001736
$ws.Range("M34").Value = -50001534
$ws.Range("H58").Value = 4713.0527
$ws.Range("I58").Value = 3058.6667
$ws.Range("J58").Value = 7549.143
$ws.Range("K58").Value = 3058.6667
$ws.Range("L58").Value = 7549.143
$ws.Range("M58").Value = -2855.6667
$ws.Range("N58").Value = -7955.143
$ws.Range("H132").Value = 5578.875
$ws.Range("I132").Value = 4705.8184
$ws.Range("J132").Value = 7499.6
$ws.Range("K132").Value = 14117.4552
$ws.Range("L132").Value = 22498.8
$ws.Range("M132").Value = -11587.4552
$ws.Range("N132").Value = -27558.8
$ws.Range("H136").Value = 4713.0527
$ws.Range("I136").Value = 3058.6667
$ws.Range("J136").Value = 7549.143
$ws.Range("K136").Value = 9176.000100000001
$ws.Range("L136").Value = 22647.429
$ws.Range("M136").Value = -6626.000100000001
$ws.Range("N136").Value = -27747.429

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1486.125
$ws.Range("I39").Value = 1199
$ws.Range("K39").Value = 3597
$ws.Range("M39").Value = -3303
$ws.Range("H122").Value = 12730.8
$ws.Range("I122").Value = 5913.5
$ws.Range("K122").Value = 53221.5
$ws.Range("M122").Value = -50771.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 10000
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H32").Value = 62483.332
$ws.Range("J32").Value = 62483.332
$ws.Range("L32").Value = 62483.332
$ws.Range("N32").Value = -63075.332
$ws.Range("H132").Value = 1948.1714
$ws.Range("I132").Value = 1012.37036
$ws.Range("J132").Value = 5106.5
$ws.Range("K132").Value = 3037.11108
$ws.Range("L132").Value = 15319.5
$ws.Range("M132").Value = -507.1110800000001
$ws.Range("N132").Value = -20379.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11880.5
$ws.Range("I7").Value = 11434.857
$ws.Range("J7").Value = 15000
$ws.Range("K7").Value = 11434.857
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = -11322.857
$ws.Range("N7").Value = -15224
$ws.Range("H43").Value = 44999.5
$ws.Range("J43").Value = 44999.5
$ws.Range("L43").Value = 44999.5
$ws.Range("N43").Value = -45385.5
$ws.Range("H46").Value = 7267.222
$ws.Range("I46").Value = 7512.1313
$ws.Range("K46").Value = 7512.1313
$ws.Range("M46").Value = -7324.1313
$ws.Range("H82").Value = 48484.332
$ws.Range("I82").Value = 885.3
$ws.Range("K82").Value = 885.3
$ws.Range("M82").Value = -524.3
$ws.Range("H85").Value = 48484.332
$ws.Range("I85").Value = 885.3
$ws.Range("K85").Value = 885.3
$ws.Range("M85").Value = 362.7
$ws.Range("H126").Value = 11880.5
$ws.Range("I126").Value = 11434.857
$ws.Range("J126").Value = 15000
$ws.Range("K126").Value = 34304.571
$ws.Range("L126").Value = 45000
$ws.Range("M126").Value = -31834.571
$ws.Range("N126").Value = -49940
$ws.Range("H136").Value = 4821.486
$ws.Range("I136").Value = 4992
$ws.Range("K136").Value = 14976
$ws.Range("M136").Value = -12426

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 17830
$ws.Range("J53").Value = 30000
$ws.Range("L53").Value = 30000
$ws.Range("N53").Value = -31214
$ws.Range("H81").Value = 5906.5713
$ws.Range("I81").Value = 3998.75
$ws.Range("J81").Value = 8450.333000000001
$ws.Range("K81").Value = 7997.5
$ws.Range("L81").Value = 16900.666
$ws.Range("M81").Value = -6936.5
$ws.Range("N81").Value = -19022.666
$ws.Range("H84").Value = 5906.5713
$ws.Range("I84").Value = 3998.75
$ws.Range("J84").Value = 8450.333000000001
$ws.Range("K84").Value = 39987.5
$ws.Range("L84").Value = 84503.33
$ws.Range("M84").Value = -34683.5
$ws.Range("N84").Value = -95111.33
$ws.Range("H96").Value = 3055.9524
$ws.Range("I96").Value = 2913.7144
$ws.Range("K96").Value = 2913.7144
$ws.Range("M96").Value = -1540.7144
$ws.Range("H100").Value = 1402.2941
$ws.Range("I100").Value = 1134
$ws.Range("K100").Value = 2268
$ws.Range("M100").Value = -1727
$ws.Range("H132").Value = 3413.34
$ws.Range("I132").Value = 3140.878
$ws.Range("K132").Value = 9422.634
$ws.Range("M132").Value = -6892.634
$ws.Range("H136").Value = 1907.5084
$ws.Range("I136").Value = 1186.4419
$ws.Range("J136").Value = 3845.375
$ws.Range("K136").Value = 3559.3257
$ws.Range("L136").Value = 11536.125
$ws.Range("M136").Value = -1009.3257
$ws.Range("N136").Value = -16636.125
